$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.902.70"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.01%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.544.95"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.36%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.34%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "205.56"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.483"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.47%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.35%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "21.29"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.54%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0581"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.81%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.78%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.765.04"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.28%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.545.56"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.40%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.27%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.78%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.885.61"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.02%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.54"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.29%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "213.40"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0683"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.19"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.78%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.36%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.97%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.16"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.54%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -3.29%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.20"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.28%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.63"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.81"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.10%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.36%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.32%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.10%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.21"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.81%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.363.32"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.71%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.95"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.52%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.975"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +6.06%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.24%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.24%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.45%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.803"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.40%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.55%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.45"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.97%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.59%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "63.28"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.09%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.38%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.678.90"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.35%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "86.25"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.65%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.79%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₇0969"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.52%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.49%  "
